$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update column F "想去人数" values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1013
$wsExhibition.Range("F5").Value = 2734
$wsExhibition.Range("F6").Value = 91
$wsExhibition.Range("F7").Value = 211
$wsExhibition.Range("F9").Value = 115
$wsExhibition.Range("F10").Value = 56
$wsExhibition.Range("F11").Value = 29
$wsExhibition.Range("F12").Value = 2548
$wsExhibition.Range("F13").Value = 680

# Sheet "全部类型" (All types) - update column F "想去人数" values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1013
$wsAll.Range("F6").Value = 2734
$wsAll.Range("F7").Value = 91
$wsAll.Range("F8").Value = 211
$wsAll.Range("F11").Value = 115
$wsAll.Range("F12").Value = 56
$wsAll.Range("F13").Value = 29
$wsAll.Range("F14").Value = 2548
$wsAll.Range("F15").Value = 680
